$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2090.7576
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 2149.8333
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 6449.499899999999
$ws.Range("M17").Value = -4332
$ws.Range("N17").Value = -6785.499899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 13754.75
$ws.Range("I21").Value = 13754.75
$ws.Range("K21").Value = 13754.75
$ws.Range("M21").Value = -13286.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 13754.75
$ws.Range("I23").Value = 13754.75
$ws.Range("K23").Value = 13754.75
$ws.Range("M23").Value = -13520.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1702.9259
$ws.Range("I41").Value = 1671.4
$ws.Range("K41").Value = 1671.4
$ws.Range("M41").Value = -1231.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 4313.4614
$ws.Range("J80").Value = 5102.5
$ws.Range("L80").Value = 15307.5
$ws.Range("N80").Value = -17303.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 4313.4614
$ws.Range("J83").Value = 5102.5
$ws.Range("L83").Value = 45922.5
$ws.Range("N83").Value = -55906.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6597.0527
$ws.Range("I113").Value = 4147.3335
$ws.Range("K113").Value = 4147.3335
$ws.Range("M113").Value = -893.3334999999997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 543.17645
$ws.Range("I115").Value = 564.625
$ws.Range("K115").Value = 1693.875
$ws.Range("M115").Value = -126.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 92388.07000000001
$ws.Range("J133").Value = 92388.07000000001
$ws.Range("L133").Value = 92388.07000000001
$ws.Range("N133").Value = -102508.07

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 59718.25
$ws.Range("J134").Value = 68991
$ws.Range("L134").Value = 68991
$ws.Range("N134").Value = -79131

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 10250
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 10250
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 10250
$ws.Range("N10").Value = -10590
$ws.Range("M10").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 13366667
$ws.Range("I11").Value = 13366667
$ws.Range("K11").Value = 13366667
$ws.Range("M11").Value = -13366523

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12223.12
$ws.Range("I32").Value = 7650.8
$ws.Range("J32").Value = 22891.867
$ws.Range("K32").Value = 7650.8
$ws.Range("L32").Value = 22891.867
$ws.Range("M32").Value = -7363.8
$ws.Range("N32").Value = -23465.867

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1999.5
$ws.Range("I63").Value = 1999.5
$ws.Range("K63").Value = 1999.5
$ws.Range("M63").Value = -1313.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1999.5
$ws.Range("I66").Value = 1999.5
$ws.Range("K66").Value = 9997.5
$ws.Range("M66").Value = -6565.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1774.0834
$ws.Range("J88").Value = 1821.2222
$ws.Range("L88").Value = 1821.2222
$ws.Range("N88").Value = -2633.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1774.0834
$ws.Range("J91").Value = 1821.2222
$ws.Range("L91").Value = 1821.2222
$ws.Range("N91").Value = -4629.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2253
$ws.Range("I132").Value = 1908.3334
$ws.Range("K132").Value = 5725.0002
$ws.Range("M132").Value = -3195.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2363
$ws.Range("I20").Value = 1913.7646
$ws.Range("K20").Value = 1913.7646
$ws.Range("M20").Value = -1666.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 205799.6
$ws.Range("I6").Value = 999
$ws.Range("K6").Value = 999
$ws.Range("M6").Value = -886

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9347.717000000001
$ws.Range("I31").Value = 1844.3636
$ws.Range("J31").Value = 21728.25
$ws.Range("K31").Value = 1844.3636
$ws.Range("L31").Value = 21728.25
$ws.Range("M31").Value = -1549.3636
$ws.Range("N31").Value = -22318.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9347.717000000001
$ws.Range("I34").Value = 1844.3636
$ws.Range("J34").Value = 21728.25
$ws.Range("K34").Value = 1844.3636
$ws.Range("L34").Value = 21728.25
$ws.Range("M34").Value = -1642.3636
$ws.Range("N34").Value = -22132.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 50001
$ws.Range("J36").Value = 50001
$ws.Range("L36").Value = 50001
$ws.Range("N36").Value = -50777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 50001
$ws.Range("J40").Value = 50001
$ws.Range("L40").Value = 50001
$ws.Range("N40").Value = -50321

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 11000
$ws.Range("J55").Value = 11000
$ws.Range("L55").Value = 11000
$ws.Range("N55").Value = -11630

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3012.625
$ws.Range("I62").Value = 2871.5715
$ws.Range("K62").Value = 2871.5715
$ws.Range("M62").Value = -2247.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3012.625
$ws.Range("I65").Value = 2871.5715
$ws.Range("K65").Value = 14357.8575
$ws.Range("M65").Value = -11237.8575

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2606695.8
$ws.Range("I99").Value = 2088.75
$ws.Range("K99").Value = 2088.75
$ws.Range("M99").Value = -590.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2606695.8
$ws.Range("I126").Value = 2088.75
$ws.Range("K126").Value = 6266.25
$ws.Range("M126").Value = -3796.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 17166.666
$ws.Range("J129").Value = 17166.666
$ws.Range("L129").Value = 17166.666
$ws.Range("N129").Value = -27166.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 64897.8
$ws.Range("J135").Value = 64897.8
$ws.Range("L135").Value = 64897.8
$ws.Range("N135").Value = -75037.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 37464
$ws.Range("I131").Value = 77866.766
$ws.Range("J131").Value = 2448.2666
$ws.Range("K131").Value = 233600.298
$ws.Range("L131").Value = 7344.7998
$ws.Range("M131").Value = -228560.298
$ws.Range("N131").Value = -17424.7998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5592
$ws.Range("I133").Value = 4122.8335
$ws.Range("J133").Value = 9999.5
$ws.Range("K133").Value = 12368.5005
$ws.Range("L133").Value = 29998.5
$ws.Range("M133").Value = -7308.500499999998
$ws.Range("N133").Value = -40118.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2119.35
$ws.Range("J137").Value = 2662.9
$ws.Range("L137").Value = 7988.700000000001
$ws.Range("N137").Value = -18188.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3722.2593
$ws.Range("I138").Value = 1812.5625
$ws.Range("K138").Value = 5437.6875
$ws.Range("M138").Value = -297.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3957871
$ws.Range("I11").Value = 3508130.5
$ws.Range("J11").Value = 5250875
$ws.Range("K11").Value = 3508130.5
$ws.Range("L11").Value = 5250875
$ws.Range("M11").Value = -3507991.5
$ws.Range("N11").Value = -5251153

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 2444.4443
$ws.Range("I12").Value = 900
$ws.Range("J12").Value = 2637.5
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 2637.5
$ws.Range("M12").Value = -760
$ws.Range("N12").Value = -2917.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 18310
$ws.Range("I52").Value = 10030
$ws.Range("J52").Value = 22450
$ws.Range("K52").Value = 10030
$ws.Range("L52").Value = 22450
$ws.Range("M52").Value = -9771
$ws.Range("N52").Value = -22968

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1855.3334
$ws.Range("I132").Value = 1321.2941
$ws.Range("K132").Value = 3963.8823
$ws.Range("M132").Value = -1433.8823

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 42505
$ws.Range("I5").Value = 35009
$ws.Range("J5").Value = 50001
$ws.Range("K5").Value = 35009
$ws.Range("L5").Value = 50001
$ws.Range("M5").Value = -34896
$ws.Range("N5").Value = -50227

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 25000
$ws.Range("J14").Value = 25000
$ws.Range("L14").Value = 25000
$ws.Range("N14").Value = -25344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 2130.6667
$ws.Range("J31").Value = 3047.75
$ws.Range("L31").Value = 3047.75
$ws.Range("N31").Value = -3543.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 41
$ws.Range("I48").Value = 41
$ws.Range("K48").Value = 41
$ws.Range("M48").Value = 620

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 2500450
$ws.Range("I17").Value = 2500450
$ws.Range("K17").Value = 2500450
$ws.Range("M17").Value = -2500278

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 25000
$ws.Range("J18").Value = 25000
$ws.Range("L18").Value = 25000
$ws.Range("N18").Value = -25346

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 43511
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 43511
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 43511
$ws.Range("N20").Value = -43991
$ws.Range("M20").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2149
$ws.Range("I81").Value = 1299
$ws.Range("J81").Value = 2999
$ws.Range("K81").Value = 2598
$ws.Range("L81").Value = 5998
$ws.Range("M81").Value = -1537
$ws.Range("N81").Value = -8120

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2149
$ws.Range("I84").Value = 1299
$ws.Range("J84").Value = 2999
$ws.Range("K84").Value = 12990
$ws.Range("L84").Value = 29990
$ws.Range("M84").Value = -7686
$ws.Range("N84").Value = -40598

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 895766.5
$ws.Range("I132").Value = 11977.033
$ws.Range("J132").Value = 2291223.5
$ws.Range("K132").Value = 35931.099
$ws.Range("L132").Value = 6873670.5
$ws.Range("M132").Value = -33401.099
$ws.Range("N132").Value = -6878730.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1813.196
$ws.Range("I136").Value = 928.60974
$ws.Range("K136").Value = 2785.82922
$ws.Range("M136").Value = -235.8292200000001
